$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sales Entries")

# New header cell L1 - "Custom Sales Ledger Name"
$ws.Range("L1").Value = "Custom Sales Ledger Name"

# Copy the style of K1 (last header) onto L1 so fill/font/border match, then
# adjust the number format to Text (matches style xf numFmtId 49).
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("L1").NumberFormat = "@"

# Column widths: K shrinks to fit its new (shorter) content, L sized to fit
# the new header text.
$ws.Columns.Item(11).ColumnWidth = 14.5703125
$ws.Columns.Item(12).ColumnWidth = 25.28515625

# Scroll/selection state recorded in the sheetView
$ws.Application.ActiveWindow.ScrollColumn = 7
$ws.Range("L1:L1048576").Select()
